# lab_1 bugfix, report update, tests added
#
# Updates the "Количество сравнений" (comparison-count) figures for the
# good_t_3 test case in the second results table (rows 13-22), and moves
# the active selection to reflect where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E19").Value = 80
$ws.Range("E20").Value = 107
$ws.Range("E21").Value = 434
$ws.Range("E22").Value = 510

$ws.Range("I23").Select()
